# Update countries & provincias Spain
#
# The "Pais" sheet is sorted descending by "Casos totales" (column B).
# Two countries (Bolivia and Haiti) received updated case counts, which
# pushed them up past the countries immediately above them in the ranking.
# That re-sort is expressed here as: give the moved-up country its fresh
# numbers on the row it now occupies, and shift the row(s) it jumped over
# down by one (re-using the numbers they already had).  The trailing
# "Datos actualizados ..." timestamp is also bumped forward by 30 minutes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bolivia overtakes Ghana (rows 65-66) ---------------------------------
$ws.Range("A65").Value = "Bolivia"
$ws.Range("B65").Value = 7136
$ws.Range("C65").Value = 476
$ws.Range("D65").Value = 677
$ws.Range("E65").Value = 6185
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 13
$ws.Range("H65").Value = 274

$ws.Range("A66").Value = "Ghana"
$ws.Range("B66").Value = 7117
$ws.Range("C66").Value = 0
$ws.Range("D66").Value = 2317
$ws.Range("E66").Value = 4766
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 34

# --- Haiti overtakes Libano, Mali and Hong Kong (rows 106-109) -----------
$ws.Range("A106").Value = "Haiti"
$ws.Range("B106").Value = 1174
$ws.Range("C106").Value = 111
$ws.Range("D106").Value = 22
$ws.Range("E106").Value = 1119
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 33

$ws.Range("A107").Value = "Libano"
$ws.Range("B107").Value = 1140
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 689
$ws.Range("E107").Value = 425
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 26

$ws.Range("A108").Value = "Mali"
$ws.Range("B108").Value = 1077
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 617
$ws.Range("E108").Value = 390
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 70

$ws.Range("A109").Value = "Hong Kong"
$ws.Range("B109").Value = 1066
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 1033
$ws.Range("E109").Value = 29
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 4

# --- Bump the "last updated" timestamp ------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 04:35"
